# Apply the "feat: Added docstrings and README." commit's spreadsheet edits.
# (The commit message is generic; the actual change updates scenario numbers
# across the Data / Pricing Table / Nutriant Table / Discounts to Volume
# Table sheets, and shifts which sheet/cell is active.)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Data" sheet (rows 2-193): columns F (Waste Diversion TPA) and
#    G (Soil Sequestration TPA) get new, rounder scenario numbers.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Data")

# Column F repeats a 12-row cycle of six baseline figures (each appearing on
# two consecutive rows). Two of the six were entered as formulas in the
# original workbook (e.g. "=309598") and stay formulas with the new literal;
# the rest are plain numbers.
$fCycleValues     = @(295000, 295000, 150000, 150000, 3000,  3000,  210000, 210000, 25000, 25000, 2500,  2500)
$fCycleIsFormula  = @($true,  $true,  $false, $false, $true, $true, $false, $false, $false,$false,$false,$false)

# Column G is constant across blocks of rows (12 rows, occasionally 24), and
# every one of those block values is now written as a plain number - any
# cell that used to hold a formula (e.g. "=17417+31519") loses the formula.
$gBlocks = @(
    @{Start=2;   End=25;  Value=15000}
    @{Start=26;  End=37;  Value=45000}
    @{Start=38;  End=49;  Value=30000}
    @{Start=50;  End=73;  Value=35000}
    @{Start=74;  End=85;  Value=65000}
    @{Start=86;  End=97;  Value=50000}
    @{Start=98;  End=121; Value=55000}
    @{Start=122; End=133; Value=85000}
    @{Start=134; End=145; Value=70000}
    @{Start=146; End=169; Value=65000}
    @{Start=170; End=181; Value=90000}
    @{Start=182; End=193; Value=75000}
)

for ($r = 2; $r -le 193; $r++) {
    $idx  = ($r - 2) % 12
    $fval = $fCycleValues[$idx]
    if ($fCycleIsFormula[$idx]) {
        $ws2.Range("F$r").Formula = "=$fval"
    } else {
        $ws2.Range("F$r").Value = $fval
    }
}

foreach ($blk in $gBlocks) {
    for ($r = $blk.Start; $r -le $blk.End; $r++) {
        $ws2.Range("G$r").Value = $blk.Value
    }
}

# ---------------------------------------------------------------------------
# 2. "Pricing Table" sheet - row 2 scenario inputs.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Pricing Table")
$ws3.Range("A2").Value = 2.25
$ws3.Range("B2").Value = 4.25
$ws3.Range("C2").Value = 1.75
$ws3.Range("D2").Value = 275000000
$ws3.Range("E2").Value = 100000000
$ws3.Range("F2").Value = 5.15
$ws3.Range("G2").Value = 3.75
$ws3.Range("H2").Value = 0.15
$ws3.Range("I2").Value = 1.05

# ---------------------------------------------------------------------------
# 3. "Nutriant Table" sheet - rows 2-5 recompute with new discount inputs.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Nutriant Table")

$ws4.Range("B2").Value = 0.19750000000000001
$ws4.Range("E2").Value = 0.249
$ws4.Range("F2").Value = 0.03
$ws4.Range("G2").Value = 0.9
$ws4.Range("H2").Value = 1.7136363636363601
$ws4.Range("I2").Value = 1.30833333333333
$ws4.Range("J2").Value = 0.78676470588235303
$ws4.Range("K2").Value = 0.53369565217391302
$ws4.Range("L2").Value = 0.72343749999999996

$ws4.Range("B3").Value = 0.185
$ws4.Range("E3").Value = 0.23599999999999999
$ws4.Range("F3").Value = 0.08
$ws4.Range("G3").Value = 0.85
$ws4.Range("H3").Value = 1.6181818181818199
$ws4.Range("I3").Value = 1.25
$ws4.Range("J3").Value = 0.65588235294117603
$ws4.Range("K3").Value = 0.46086956521739098
$ws4.Range("L3").Value = 0.69062500000000004

$ws4.Range("B4").Value = 0.17249999999999999
$ws4.Range("E4").Value = 0.184
$ws4.Range("F4").Value = 0.12
$ws4.Range("G4").Value = 0.8
$ws4.Range("H4").Value = 1.52272727272727
$ws4.Range("I4").Value = 0.89166666666666705
$ws4.Range("J4").Value = 0.42499999999999999
$ws4.Range("K4").Value = 0.33804347826087
$ws4.Range("L4").Value = 0.45781250000000001

$ws4.Range("B5").Value = 0.155
$ws4.Range("E5").Value = 0.17199999999999999
$ws4.Range("F5").Value = 0.22
$ws4.Range("G5").Value = 0.7
$ws4.Range("H5").Value = 1.4272727272727299
$ws4.Range("I5").Value = 0.73333333333333295
$ws4.Range("J5").Value = 0.39411764705882402
$ws4.Range("K5").Value = 0.26521739130434802
$ws4.Range("L5").Value = 0.42499999999999999

# ---------------------------------------------------------------------------
# 4. "Discounts to Volume Table" sheet - Wood/Land Use discount cells.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Discounts to Volume Table")
$ws5.Range("C2").Value = 0.75
$ws5.Range("D2").Value = 0.5
$ws5.Range("C5").Value = 0.85
$ws5.Range("D5").Value = 0.55000000000000004
$ws5.Range("C6").Value = 0.25
$ws5.Range("D6").Value = 0.15
$ws5.Range("C8").Value = 0.6
$ws5.Range("D8").Value = 0.25
$ws5.Range("C9").Value = 0.4
$ws5.Range("D9").Value = 0.25

# ---------------------------------------------------------------------------
# 5. View state: active sheet moves from "Scenario Input" to
#    "Discounts to Volume Table" (index 4), and each sheet remembers a new
#    selection; the "Data" sheet also scrolls so row 166 is near the top.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Scenario Input")

$ws1.Activate()
$ws1.Range("A5").Select()

$ws2.Activate()
$ws2.Range("G186").Select()
$excel.ActiveWindow.ScrollRow = 166

$ws3.Activate()
$ws3.Range("I2").Select()

$ws4.Activate()
$ws4.Range("L6").Select()

$ws5.Activate()
$ws5.Range("D10").Select()
